$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 895.43335
$ws.Range("I46").Value = 200
$ws.Range("J46").Value = 919.4138
$ws.Range("K46").Value = 600
$ws.Range("L46").Value = 2758.2414
$ws.Range("M46").Value = -481
$ws.Range("N46").Value = -2996.2414
$ws.Range("H60").Value = 895.43335
$ws.Range("I60").Value = 200
$ws.Range("J60").Value = 919.4138
$ws.Range("K60").Value = 600
$ws.Range("L60").Value = 2758.2414
$ws.Range("M60").Value = -116
$ws.Range("N60").Value = -3726.2414
$ws.Range("H64").Value = 58081.61
$ws.Range("J64").Value = 2691.3333
$ws.Range("L64").Value = 2691.3333
$ws.Range("N64").Value = -3187.3333
$ws.Range("H67").Value = 58081.61
$ws.Range("J67").Value = 2691.3333
$ws.Range("L67").Value = 2691.3333
$ws.Range("N67").Value = -4407.3333
$ws.Range("H94").Value = 100001150
$ws.Range("I94").Value = 1280.4445
$ws.Range("K94").Value = 1280.4445
$ws.Range("M94").Value = -829.4445000000001
$ws.Range("H132").Value = 4468931
$ws.Range("I132").Value = 4812289
$ws.Range("K132").Value = 14436867
$ws.Range("M132").Value = -14434337
$ws.Range("H137").Value = 1290.1372
$ws.Range("I137").Value = 1097.6511
$ws.Range("J137").Value = 2324.75
$ws.Range("K137").Value = 3292.9533
$ws.Range("L137").Value = 6974.25
$ws.Range("M137").Value = -742.9533000000001
$ws.Range("N137").Value = -12074.25
$ws.Range("H138").Value = 4908.712
$ws.Range("I138").Value = 972.26
$ws.Range("J138").Value = 26777.889
$ws.Range("K138").Value = 2916.78
$ws.Range("L138").Value = 80333.667
$ws.Range("M138").Value = 2223.22
$ws.Range("N138").Value = -90613.667
$ws.Range("H141").Value = 2501.4595
$ws.Range("I141").Value = 2367.1562
$ws.Range("K141").Value = 7101.4686
$ws.Range("M141").Value = -1921.4686

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 22289.906
$ws.Range("I32").Value = 4655.69
$ws.Range("K32").Value = 4655.69
$ws.Range("M32").Value = -4368.69
$ws.Range("H61").Value = 2106.7856
$ws.Range("I61").Value = 1043.3889
$ws.Range("K61").Value = 1043.3889
$ws.Range("M61").Value = -831.3888999999999
$ws.Range("H110").Value = 167000660
$ws.Range("J110").Value = 1304.3334
$ws.Range("L110").Value = 1304.3334
$ws.Range("N110").Value = -5394.3334
$ws.Range("H122").Value = 1366.1389
$ws.Range("I122").Value = 1340.6177
$ws.Range("K122").Value = 4021.8531
$ws.Range("M122").Value = -1571.8531
$ws.Range("H133").Value = 38574.5
$ws.Range("J133").Value = 38574.5
$ws.Range("L133").Value = 38574.5
$ws.Range("N133").Value = -43634.5
$ws.Range("H136").Value = 2106.7856
$ws.Range("I136").Value = 1043.3889
$ws.Range("K136").Value = 3130.1667
$ws.Range("M136").Value = -580.1666999999998

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 648.1667
$ws.Range("I64").Value = 176.5
$ws.Range("J64").Value = 884
$ws.Range("K64").Value = 176.5
$ws.Range("L64").Value = 884
$ws.Range("M64").Value = 48.5
$ws.Range("N64").Value = -1334
$ws.Range("H67").Value = 648.1667
$ws.Range("I67").Value = 176.5
$ws.Range("J67").Value = 884
$ws.Range("K67").Value = 176.5
$ws.Range("L67").Value = 884
$ws.Range("M67").Value = 603.5
$ws.Range("N67").Value = -2444
$ws.Range("H99").Value = 1582.119
$ws.Range("I99").Value = 1104.45
$ws.Range("J99").Value = 2016.3636
$ws.Range("K99").Value = 1104.45
$ws.Range("L99").Value = 2016.3636
$ws.Range("M99").Value = 393.55
$ws.Range("N99").Value = -5012.3636
$ws.Range("H134").Value = 3183.3125
$ws.Range("I134").Value = 3014.9583
$ws.Range("J134").Value = 3688.375
$ws.Range("K134").Value = 9044.874899999999
$ws.Range("L134").Value = 11065.125
$ws.Range("M134").Value = -6509.874899999999
$ws.Range("N134").Value = -16135.125

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 980.3
$ws.Range("I16").Value = 972.8570999999999
$ws.Range("K16").Value = 972.8570999999999
$ws.Range("M16").Value = -685.8570999999999
$ws.Range("H31").Value = 24481.428
$ws.Range("I31").Value = 1201.08
$ws.Range("K31").Value = 1201.08
$ws.Range("M31").Value = -906.0799999999999
$ws.Range("H34").Value = 24481.428
$ws.Range("I34").Value = 1201.08
$ws.Range("K34").Value = 1201.08
$ws.Range("M34").Value = -999.0799999999999
$ws.Range("H113").Value = 980.3
$ws.Range("I113").Value = 972.8570999999999
$ws.Range("K113").Value = 972.8570999999999
$ws.Range("M113").Value = 1197.1429

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 400
$ws.Range("J17").Value = 400
$ws.Range("L17").Value = 1200
$ws.Range("N17").Value = -1538
$ws.Range("H18").Value = 355.33334
$ws.Range("I18").Value = 322.86667
$ws.Range("J18").Value = 517.6667
$ws.Range("K18").Value = 968.60001
$ws.Range("L18").Value = 1553.0001
$ws.Range("M18").Value = -799.60001
$ws.Range("N18").Value = -1891.0001
$ws.Range("H32").Value = 6687000
$ws.Range("I32").Value = 20001000
$ws.Range("J32").Value = 30000
$ws.Range("K32").Value = 60003000
$ws.Range("L32").Value = 90000
$ws.Range("M32").Value = -60002717
$ws.Range("N32").Value = -90566
$ws.Range("H34").Value = 694.7143
$ws.Range("J34").Value = 944.125
$ws.Range("L34").Value = 2832.375
$ws.Range("N34").Value = -3000.375
$ws.Range("H39").Value = 2133.3333
$ws.Range("J39").Value = 2950
$ws.Range("L39").Value = 8850
$ws.Range("N39").Value = -9438
$ws.Range("H55").Value = 16483.334
$ws.Range("J55").Value = 7858.3335
$ws.Range("L55").Value = 23575.0005
$ws.Range("N55").Value = -23929.0005
$ws.Range("H131").Value = 7510.9536
$ws.Range("J131").Value = 8272.103999999999
$ws.Range("L131").Value = 24816.312
$ws.Range("N131").Value = -34896.312
$ws.Range("H140").Value = 5770.48
$ws.Range("I140").Value = 10090.417
$ws.Range("J140").Value = 1782.8462
$ws.Range("K140").Value = 30271.251
$ws.Range("L140").Value = 5348.5386
$ws.Range("M140").Value = -25091.251
$ws.Range("N140").Value = -15708.5386

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 1156814.9
$ws.Range("I43").Value = 1875449.2
$ws.Range("K43").Value = 1875449.2
$ws.Range("M43").Value = -1875298.2
$ws.Range("H70").Value = 104896.1
$ws.Range("I70").Value = 203791.6
$ws.Range("J70").Value = 6000.6
$ws.Range("K70").Value = 203791.6
$ws.Range("L70").Value = 6000.6
$ws.Range("M70").Value = -203521.6
$ws.Range("N70").Value = -6540.6
$ws.Range("H73").Value = 104896.1
$ws.Range("I73").Value = 203791.6
$ws.Range("J73").Value = 6000.6
$ws.Range("K73").Value = 203791.6
$ws.Range("L73").Value = 6000.6
$ws.Range("M73").Value = -202855.6
$ws.Range("N73").Value = -7872.6
$ws.Range("H123").Value = 7018.5713
$ws.Range("J123").Value = 7018.5713
$ws.Range("L123").Value = 7018.5713
$ws.Range("N123").Value = -11918.5713
$ws.Range("H132").Value = 2403.0244
$ws.Range("I132").Value = 1655.1072
$ws.Range("K132").Value = 4965.321599999999
$ws.Range("M132").Value = -2435.321599999999
$ws.Range("H135").Value = 27800
$ws.Range("J135").Value = 27800
$ws.Range("L135").Value = 27800
$ws.Range("N135").Value = -37940

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 632982.4399999999
$ws.Range("I46").Value = 188.375
$ws.Range("J46").Value = 1265776.5
$ws.Range("K46").Value = 188.375
$ws.Range("L46").Value = 1265776.5
$ws.Range("M46").Value = -0.375
$ws.Range("N46").Value = -1266152.5
$ws.Range("H132").Value = 2512.5112
$ws.Range("I132").Value = 2535.1843
$ws.Range("J132").Value = 2389.4285
$ws.Range("K132").Value = 7605.5529
$ws.Range("L132").Value = 7168.2855
$ws.Range("M132").Value = -5075.5529
$ws.Range("N132").Value = -12228.2855
$ws.Range("H133").Value = 46800
$ws.Range("J133").Value = 46800
$ws.Range("L133").Value = 46800
$ws.Range("N133").Value = -51860
$ws.Range("H136").Value = 1391.303
$ws.Range("I136").Value = 1191.5358
$ws.Range("J136").Value = 2510
$ws.Range("K136").Value = 3574.6074
$ws.Range("L136").Value = 7530
$ws.Range("M136").Value = -1024.6074
$ws.Range("N136").Value = -12630

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H16").Value = 36000
$ws.Range("J16").Value = 36000
$ws.Range("L16").Value = 36000
$ws.Range("N16").Value = -36584
$ws.Range("H122").Value = 1926.3684
$ws.Range("I122").Value = 1601
$ws.Range("J122").Value = 1987.375
$ws.Range("K122").Value = 4803
$ws.Range("L122").Value = 5962.125
$ws.Range("M122").Value = -2353
$ws.Range("N122").Value = -10862.125
$ws.Range("H132").Value = 1659.2693
$ws.Range("I132").Value = 1635.3175
$ws.Range("K132").Value = 4905.9525
$ws.Range("L132").Value = 5279.6001
$ws.Range("M132").Value = -2375.9525
$ws.Range("N132").Value = -10339.6001
$ws.Range("H136").Value = 713.8043
$ws.Range("I136").Value = 415
$ws.Range("J136").Value = 1560.4166
$ws.Range("K136").Value = 1245
$ws.Range("L136").Value = 4681.2498
$ws.Range("M136").Value = 1305
$ws.Range("N136").Value = -9781.2498
